$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates (row 2) ---
$ws.Range("P2").Value = "request_Id=184"
$ws.Range("K2").Value = "select iqr.id, iq.insurance_premium_amount, iq.insurance_premium_currency, iq.policy_limit_amount from insurancequotes iq INNER JOIN insurancequoterequests iqr on iq.id = iqr.insurance_quote_id and iqr.id  =  [request_Id]"
$ws.Range("M2").Value = "policy_limit_amount=[0].policy_limit_amount;quoteId=[0].id;"

# --- Cell content updates (row 3) ---
$ws.Range("N3").Value = "select iqr.id, iq.insurance_premium_amount, iq.insurance_premium_currency, iq.policy_limit_amount from insurancequotes iq INNER JOIN insurancequoterequests iqr on iq.id = iqr.insurance_quote_id and iqr.id  =  [request_Id]`nid,insurance_premium_amount, insurance_premium_currency, policy_limit_amount`n[quoteId],d~500.00,CHF,d~50000.00"

# --- Column M width change ---
$ws.Columns.Item(13).ColumnWidth = 59.75

# --- Sheet view changes: scroll position and selection ---
$ws.Range("N3").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 12
